$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# Header row (row 1): B..K
$ws.Range("B1").Value = 'company'
$ws.Range("C1").Value = 'name'
$ws.Range("D1").Value = 'owner'
$ws.Range("E1").Value = 'property_category'
$ws.Range("F1").Value = 'category'
$ws.Range("G1").Value = 'date'
$ws.Range("H1").Value = 'legislator_name'
$ws.Range("I1").Value = 'legislator_id'
$ws.Range("J1").Value = 'source_file'
$ws.Range("K1").Value = 'index'

# Data rows 2..19
# row 2
$ws.Range("B2").Value = '富邦人壽'
$ws.Range("C2").Value = '终身壽險丙型'
$ws.Range("D2").Value = '陳秀蘭'
$ws.Range("E2").Value = 'insurance'
$ws.Range("F2").Value = 'normal'
$ws.Range("G2").Value = '2012-04-27'
$ws.Range("H2").Value = '陳雪生'
$ws.Range("I2").Value = 1751
$ws.Range("J2").Value = 'tmp5a001'
$ws.Range("K2").Value = 103

# row 3
$ws.Range("B3").Value = '富邦人壽'
$ws.Range("C3").Value = '添富萬能壽險'
$ws.Range("D3").Value = '陳秀蘭'
$ws.Range("E3").Value = 'insurance'
$ws.Range("F3").Value = 'normal'
$ws.Range("G3").Value = '2012-04-27'
$ws.Range("H3").Value = '陳雪生'
$ws.Range("I3").Value = 1751
$ws.Range("J3").Value = 'tmp5a001'
$ws.Range("K3").Value = 104

# row 4
$ws.Range("B4").Value = '富邦人壽'
$ws.Range("C4").Value = '豐帘養老保險'
$ws.Range("D4").Value = '陳秀蘭'
$ws.Range("E4").Value = 'insurance'
$ws.Range("F4").Value = 'normal'
$ws.Range("G4").Value = '2012-04-27'
$ws.Range("H4").Value = '陳雪生'
$ws.Range("I4").Value = 1751
$ws.Range("J4").Value = 'tmp5a001'
$ws.Range("K4").Value = 105

# row 5
$ws.Range("B5").Value = '富邦人壽'
$ws.Range("C5").Value = '雙福還本分紅終身保險'
$ws.Range("D5").Value = '陳秀蘭'
$ws.Range("E5").Value = 'insurance'
$ws.Range("F5").Value = 'normal'
$ws.Range("G5").Value = '2012-04-27'
$ws.Range("H5").Value = '陳雪生'
$ws.Range("I5").Value = 1751
$ws.Range("J5").Value = 'tmp5a001'
$ws.Range("K5").Value = 107

# row 6
$ws.Range("B6").Value = '富邦人壽'
$ws.Range("C6").Value = '全福還本分紅終身保險'
$ws.Range("D6").Value = '陳秀蘭'
$ws.Range("E6").Value = 'insurance'
$ws.Range("F6").Value = 'normal'
$ws.Range("G6").Value = '2012-04-27'
$ws.Range("H6").Value = '陳雪生'
$ws.Range("I6").Value = 1751
$ws.Range("J6").Value = 'tmp5a001'
$ws.Range("K6").Value = 108

# row 7
$ws.Range("B7").Value = '富邦人壽'
$ws.Range("C7").Value = '大多利率變動型年金保險(甲型）'
$ws.Range("D7").Value = '陳秀蘭'
$ws.Range("E7").Value = 'insurance'
$ws.Range("F7").Value = 'normal'
$ws.Range("G7").Value = '2012-04-27'
$ws.Range("H7").Value = '陳雪生'
$ws.Range("I7").Value = 1751
$ws.Range("J7").Value = 'tmp5a001'
$ws.Range("K7").Value = 109

# row 8
$ws.Range("B8").Value = '富邦人壽'
$ws.Range("C8").Value = '心得意利率變動型年金保險(甲型）'
$ws.Range("D8").Value = '陳秀蘭'
$ws.Range("E8").Value = 'insurance'
$ws.Range("F8").Value = 'normal'
$ws.Range("G8").Value = '2012-04-27'
$ws.Range("H8").Value = '陳雪生'
$ws.Range("I8").Value = 1751
$ws.Range("J8").Value = 'tmp5a001'
$ws.Range("K8").Value = 110

# row 9
$ws.Range("B9").Value = '富邦人壽'
$ws.Range("C9").Value = '豐華養老保險'
$ws.Range("D9").Value = '陳秀蘭'
$ws.Range("E9").Value = 'insurance'
$ws.Range("F9").Value = 'normal'
$ws.Range("G9").Value = '2012-04-27'
$ws.Range("H9").Value = '陳雪生'
$ws.Range("I9").Value = 1751
$ws.Range("J9").Value = 'tmp5a001'
$ws.Range("K9").Value = 111

# row 10
$ws.Range("B10").Value = '富邦人壽'
$ws.Range("C10").Value = '五五得利終身壽險'
$ws.Range("D10").Value = '陳秀蘭'
$ws.Range("E10").Value = 'insurance'
$ws.Range("F10").Value = 'normal'
$ws.Range("G10").Value = '2012-04-27'
$ws.Range("H10").Value = '陳雪生'
$ws.Range("I10").Value = 1751
$ws.Range("J10").Value = 'tmp5a001'
$ws.Range("K10").Value = 112

# row 11
$ws.Range("B11").Value = '富邦人壽'
$ws.Range("C11").Value = '豐華養老保險'
$ws.Range("D11").Value = '陳秀蘭'
$ws.Range("E11").Value = 'insurance'
$ws.Range("F11").Value = 'normal'
$ws.Range("G11").Value = '2012-04-27'
$ws.Range("H11").Value = '陳雪生'
$ws.Range("I11").Value = 1751
$ws.Range("J11").Value = 'tmp5a001'
$ws.Range("K11").Value = 113

# row 12
$ws.Range("B12").Value = '富邦人壽'
$ws.Range("C12").Value = '豐沛利率變動型養老保險'
$ws.Range("D12").Value = '陳秀蘭'
$ws.Range("E12").Value = 'insurance'
$ws.Range("F12").Value = 'normal'
$ws.Range("G12").Value = '2012-04-27'
$ws.Range("H12").Value = '陳雪生'
$ws.Range("I12").Value = 1751
$ws.Range("J12").Value = 'tmp5a001'
$ws.Range("K12").Value = 114

# row 13
$ws.Range("B13").Value = '富邦人壽'
$ws.Range("C13").Value = '美利成增外幣養老保險（美元）'
$ws.Range("D13").Value = '陳秀蘭'
$ws.Range("E13").Value = 'insurance'
$ws.Range("F13").Value = 'normal'
$ws.Range("G13").Value = '2012-04-27'
$ws.Range("H13").Value = '陳雪生'
$ws.Range("I13").Value = 1751
$ws.Range("J13").Value = 'tmp5a001'
$ws.Range("K13").Value = 115

# row 14
$ws.Range("B14").Value = '富邦人壽'
$ws.Range("C14").Value = '增美利外幣終身壽險(美元）'
$ws.Range("D14").Value = '陳秀蘭'
$ws.Range("E14").Value = 'insurance'
$ws.Range("F14").Value = 'normal'
$ws.Range("G14").Value = '2012-04-27'
$ws.Range("H14").Value = '陳雪生'
$ws.Range("I14").Value = 1751
$ws.Range("J14").Value = 'tmp5a001'
$ws.Range("K14").Value = 116

# row 15
$ws.Range("B15").Value = '幸福人壽'
$ws.Range("C15").Value = '新幸福人生養老保險'
$ws.Range("D15").Value = '陳秀蘭'
$ws.Range("E15").Value = 'insurance'
$ws.Range("F15").Value = 'normal'
$ws.Range("G15").Value = '2012-04-27'
$ws.Range("H15").Value = '陳雪生'
$ws.Range("I15").Value = 1751
$ws.Range("J15").Value = 'tmp5a001'
$ws.Range("K15").Value = 117

# row 16
$ws.Range("B16").Value = '全球人壽'
$ws.Range("C16").Value = '全球增額終身壽險B型'
$ws.Range("D16").Value = '陳秀蘭'
$ws.Range("E16").Value = 'insurance'
$ws.Range("F16").Value = 'normal'
$ws.Range("G16").Value = '2012-04-27'
$ws.Range("H16").Value = '陳雪生'
$ws.Range("I16").Value = 1751
$ws.Range("J16").Value = 'tmp5a001'
$ws.Range("K16").Value = 118

# row 17
$ws.Range("B17").Value = '遠雄人壽'
$ws.Range("C17").Value = '終身壽險定期給付型'
$ws.Range("D17").Value = '陳秀蘭'
$ws.Range("E17").Value = 'insurance'
$ws.Range("F17").Value = 'normal'
$ws.Range("G17").Value = '2012-04-27'
$ws.Range("H17").Value = '陳雪生'
$ws.Range("I17").Value = 1751
$ws.Range("J17").Value = 'tmp5a001'
$ws.Range("K17").Value = 119

# row 18
$ws.Range("B18").Value = '蘇黎世國際人壽'
$ws.Range("C18").Value = '终身壽險定期給付甲型'
$ws.Range("D18").Value = '陳秀蘭'
$ws.Range("E18").Value = 'insurance'
$ws.Range("F18").Value = 'normal'
$ws.Range("G18").Value = '2012-04-27'
$ws.Range("H18").Value = '陳雪生'
$ws.Range("I18").Value = 1751
$ws.Range("J18").Value = 'tmp5a001'
$ws.Range("K18").Value = 120

# row 19
$ws.Range("B19").Value = '遠雄人壽'
$ws.Range("C19").Value = '新雄多利保險甲型'
$ws.Range("D19").Value = '陳秀蘭'
$ws.Range("E19").Value = 'insurance'
$ws.Range("F19").Value = 'normal'
$ws.Range("G19").Value = '2012-04-27'
$ws.Range("H19").Value = '陳雪生'
$ws.Range("I19").Value = 1751
$ws.Range("J19").Value = 'tmp5a001'
$ws.Range("K19").Value = 121
